$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new parameter row (Half_Rocker_Width) below the existing table.
$ws.Range("A17").Value = "Half_Rocker_Width (mm)"
$ws.Range("B17").Value = 7.5
# Match the numeric formatting used by the rest of column B (2 decimals).
$ws.Range("B17").NumberFormat = "0.00"

# Reflect the new active cell/selection recorded in the saved file.
$ws.Range("F9").Select()
